$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.116994380950928
$ws.Range("B1").Value = 2.503241300582886
$ws.Range("C1").Value = 6.166455745697021
$ws.Range("D1").Value = 2.173532485961914
$ws.Range("E1").Value = 1.252100229263306
